$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.52"
$ws.Range("E2").Value = "'1.51%"
$ws.Range("G2").Value = "'4"

$ws.Range("D3").Value = "'40.00"
$ws.Range("E3").Value = "'-2.61%"
$ws.Range("G3").Value = "'4"

$ws.Range("D4").Value = "'5.115"
$ws.Range("E4").Value = "'-2.18%"
$ws.Range("G4").Value = "'4"

$ws.Range("D5").Value = "'0.07572"
$ws.Range("E5").Value = "'-1.18%"
$ws.Range("G5").Value = "'4"

$ws.Range("B6").Value = "'FTXToken"
$ws.Range("C6").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.675"
$ws.Range("E6").Value = "'3.22%"
$ws.Range("G6").Value = "'4"

$ws.Range("B7").Value = "'MXToken"
$ws.Range("C7").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9308"
$ws.Range("E7").Value = "'1.43%"
$ws.Range("G7").Value = "'4"

$ws.Range("B8").Value = "'BTSEToken"
$ws.Range("C8").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.424"
$ws.Range("E8").Value = "'-0.17%"
$ws.Range("G8").Value = "'4"

$ws.Range("B9").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1210"
$ws.Range("E9").Value = "'-2.93%"
$ws.Range("G9").Value = "'4"

$ws.Range("B10").Value = "'WazirX"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1812"
$ws.Range("E10").Value = "'-1.34%"
$ws.Range("G10").Value = "'4"

$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09089"
$ws.Range("E11").Value = "'-0.21%"
$ws.Range("G11").Value = "'4"

$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04145"
$ws.Range("E12").Value = "'-2.62%"
$ws.Range("G12").Value = "'4"

$ws.Range("B13").Value = "'BitMartToken"
$ws.Range("C13").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1053"
$ws.Range("E13").Value = "'0.09%"
$ws.Range("G13").Value = "'4"

$ws.Range("B14").Value = "'BitForexToken"
$ws.Range("C14").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001279"
$ws.Range("E14").Value = "'1.30%"
$ws.Range("G14").Value = "'4"

$ws.Range("B15").Value = "'TigerCash"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005838"
$ws.Range("E15").Value = "'-0.95%"
$ws.Range("G15").Value = "'4"

$ws.Range("B16").Value = "'UpBots"
$ws.Range("C16").Value = "'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007522"
$ws.Range("E16").Value = "'0.18%"
$ws.Range("G16").Value = "'4"

$ws.Range("B17").Value = "'LEO"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.351"
$ws.Range("E17").Value = "'-0.13%"
$ws.Range("G17").Value = "'4"

$ws.Range("B18").Value = "'GateToken"
$ws.Range("C18").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.312"
$ws.Range("E18").Value = "'-0.16%"
$ws.Range("G18").Value = "'4"

$ws.Range("D19").Value = "'0.3354"
$ws.Range("E19").Value = "'0.56%"
$ws.Range("G19").Value = "'4"

$ws.Range("D20").Value = "'7.643"
$ws.Range("E20").Value = "'6.25%"
$ws.Range("G20").Value = "'4"

$ws.Range("E21").Value = "'-2.38%"
$ws.Range("G21").Value = "'4"

$ws.Range("D22").Value = "'0.2810"
$ws.Range("E22").Value = "'-2.90%"
$ws.Range("G22").Value = "'4"

$ws.Range("D23").Value = "'0.04022"
$ws.Range("E23").Value = "'-1.11%"
$ws.Range("G23").Value = "'4"

$ws.Range("D24").Value = "'0.001265"
$ws.Range("E24").Value = "'0.42%"
$ws.Range("G24").Value = "'4"

$ws.Range("D25").Value = "'0.004049"
$ws.Range("E25").Value = "'-2.41%"
$ws.Range("G25").Value = "'4"

$ws.Range("D26").Value = "'0.0001270"
$ws.Range("E26").Value = "'-0.21%"
$ws.Range("G26").Value = "'4"

$ws.Range("G27").Value = "'4"

$ws.Range("G28").Value = "'4"

$ws.Range("G29").Value = "'4"

$ws.Range("G30").Value = "'4"

$ws.Range("G31").Value = "'4"

$ws.Range("G32").Value = "'4"

$ws.Range("G33").Value = "'4"

$ws.Range("G34").Value = "'4"

$ws.Range("G35").Value = "'4"

$ws.Range("G36").Value = "'4"

$ws.Range("G37").Value = "'4"

$ws.Range("D38").Value = "'0.02432"
$ws.Range("E38").Value = "'-1.47%"
$ws.Range("G38").Value = "'4"

$ws.Range("E39").Value = "'-2.50%"
$ws.Range("G39").Value = "'4"

$ws.Range("D40").Value = "'0.007697"
$ws.Range("E40").Value = "'-1.88%"
$ws.Range("G40").Value = "'4"

$ws.Range("D41").Value = "'0.1300"
$ws.Range("E41").Value = "'-1.04%"
$ws.Range("G41").Value = "'4"

$ws.Range("D42").Value = "'0.007613"
$ws.Range("E42").Value = "'11.52%"
$ws.Range("G42").Value = "'4"

$ws.Range("E43").Value = "'14.62%"
$ws.Range("G43").Value = "'4"

$ws.Range("D44").Value = "'0.008138"
$ws.Range("E44").Value = "'4.48%"
$ws.Range("G44").Value = "'4"

$ws.Range("D45").Value = "'0.3114"
$ws.Range("E45").Value = "'1.90%"
$ws.Range("G45").Value = "'4"

$ws.Range("D46").Value = "'0.00006592"
$ws.Range("E46").Value = "'-1.92%"
$ws.Range("G46").Value = "'4"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.22%"
$ws.Range("G47").Value = "'4"

$ws.Range("D48").Value = "'0.2691"
$ws.Range("E48").Value = "'58.40%"
$ws.Range("G48").Value = "'4"

$ws.Range("D49").Value = "'0.004201"
$ws.Range("E49").Value = "'2.46%"
$ws.Range("G49").Value = "'4"

$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.22%"
$ws.Range("G50").Value = "'4"

$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.22%"
$ws.Range("G51").Value = "'4"
